# Update the "Förändrad" (Changed) date column C for all data rows
# from serial date 45175 (2023-09-06) to 45177 (2023-09-08).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C533").Value = 45177
